$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11 "Marking": Right mark 4 -> 5, Wrong mark -1 -> -1.2
$ws.Range("B11").Value = 5
$ws.Range("C11").Value = -1.2

# Row 12 "Total": recalculated totals after the marking-scheme change
$ws.Range("B12").Value = 135
$ws.Range("C12").Value = -0.0
$ws.Range("E12").Value = "135.0/140"
